$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# N2 is stored as inline/text string (not a real date), keep it as text
$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 121173920.19
$ws.Range("P2").Value = 1997388921.86
$ws.Range("Q2").Value = 1773812797.91
$ws.Range("R2").Value = 19.4998330754
$ws.Range("S2").Value = 1352264936.42
$ws.Range("T2").Value = 1352264936.42
$ws.Range("U2").Value = 22.289110205
$ws.Range("V2").Value = 220630729.68
$ws.Range("W2").Value = 81733960.03
$ws.Range("X2").Value = 29854800.46
$ws.Range("Y2").Value = 159423923.94
$ws.Range("Z2").Value = 159804469.48
$ws.Range("AA2").Value = 38687693.25
$ws.Range("AG2").Value = 28394407.5
$ws.Range("AP2").Value = 20.0330962899
$ws.Range("AQ2").Value = 36.53770105913
$ws.Range("AR2").Value = 31.640847270079
$ws.Range("AS2").Value = 115848479.87
$ws.Range("AT2").Value = 22.555864322789
